# 337 - fixes emails in excel files and the related tests
#
# The "Email" column (G) on Sheet1 had a formatting bug: the trailing
# row-number digit was appended *after* the domain instead of being part
# of the local-part, e.g. "nemreg1es@mail.com1" instead of
# "nemreg1es1@mail.com". Rows 4-6 just need the text corrected. Rows 7-8
# (whose "email" values were actually just "111111111X@mail.com", not
# real emails at all) had the bad cell removed entirely rather than
# replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the malformed emails for rows 4-6.
$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"

# Rows 7-8 never had a valid email - drop those cells entirely.
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Leave the selection on the cell that was corrected.
$ws.Range("G4").Select()
